# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Prices in column D that look like plain numbers are entered with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# inline-string data, e.g. "563.68" not the number 563.68) rather than
# auto-converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.062.70"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.423.04"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'563.68"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'143.83"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "2.422.41"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -3.81%  "
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "'26.12"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D17").Value = "61.926.26"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "2.437.32"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "'11.29"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'324.02"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.85"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'67.17"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").Value = "'1.74"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("D26").Value = "'8.80"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").Value = "'556.97"
$ws.Range("E27").Value = "  -5.57%  "
$ws.Range("D28").Value = "2.543.35"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "0.0₃0934"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "'8.20"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "'1.40"
$ws.Range("E32").Value = "  -4.90%  "
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'4.75"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").Value = "'0.380"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").Value = "'5.48"
$ws.Range("E39").Value = "  -4.70%  "
$ws.Range("D40").Value = "'152.45"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Value = "'18.65"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "'1.81"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "'2.26"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").Value = "'147.62"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").Value = "'0.0531"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "'19.90"
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("D50").Value = "'0.0919"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("E51").Value = "  -0.48%  "
